# March 24 update 3
# Adds three new columns (renewd, PlanID, iteration) to the bldg sheet,
# and populates them ("after", 20160319, 16) for every existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (M1:O1), matching the header style already used by L1 ("status")
$ws.Range("M1").Value = "renewd"
$ws.Range("N1").Value = "PlanID"
$ws.Range("O1").Value = "iteration"
$ws.Range("L1").Copy()
$ws.Range("M1:O1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New data columns for every existing data row (rows 2 through 55)
for ($r = 2; $r -le 55; $r++) {
    $ws.Cells.Item($r, 13).Value = "after"      # column M - renewd
    $ws.Cells.Item($r, 14).Value = 20160319     # column N - PlanID
    $ws.Cells.Item($r, 15).Value = 16           # column O - iteration
}
